# Update "want to attend" counts (column F) for a few convention entries
# in both the "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 15, 26, 31
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F15").Value = 1113
$wsExhibit.Range("F26").Value = 1090
$wsExhibit.Range("F31").Value = 315

# Sheet "全部类型" - rows 16, 27, 32 (same events, shifted by one row)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F16").Value = 1113
$wsAll.Range("F27").Value = 1090
$wsAll.Range("F32").Value = 315
